$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 45: predidx / pred_name change
$ws.Range("D45").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E45").Value = "['Normal']"

# Row 88: predidx / pred_name change
$ws.Range("D88").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E88").Value = "['Normal']"

# Row 118: predidx / pred_name change
$ws.Range("D118").Value = "[1, 0, 0, 1, 1, 0, 0]"
$ws.Range("E118").Value = "['Normal', 'ParamViolation', 'RegulationViolation']"
